# Update CDA Logical model for ST.r2b
# Applies the changes described by the diff:
#  - Rename "Include from RoleCode" sheet to "Include #0"
#  - Update Version and Date values on the Metadata sheet
#  - Insert a new "Jurisdiction" property row on the Metadata sheet
#  - Fill in the (previously empty) Description value
#  - Shift Purpose / Copyright / Immutable rows down by one row

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsInclude = $wb.Worksheets.Item("Include from RoleCode")

# 1. Rename the include sheet
$wsInclude.Name = "Include #0"

# 2. Update simple Version / Date values
$wsMeta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$wsMeta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# 3. Make room for the new "Jurisdiction" row by shifting rows 11-14
#    (Description, Purpose, Copyright, Immutable) down to rows 12-15.
#    Work from the bottom up so we don't overwrite data we still need.

# Prepare row 15 (new) with the same formatting as row 14, then fill in
# the values that used to live in row 14 (Immutable | BooleanType[null]).
$wsMeta.Range("A14:B14").Copy()
$wsMeta.Range("A15:B15").PasteSpecial(-4122)
$wsMeta.Range("A15").Value = $wsMeta.Range("A14").Value()
$wsMeta.Range("B15").Value = $wsMeta.Range("B14").Value()

# Row 14 now becomes the old row 13 (Copyright | "")
$wsMeta.Range("A14").Value = $wsMeta.Range("A13").Value()
$wsMeta.Range("B14").Value = ""

# Row 13 now becomes the old row 12 (Purpose | "")
$wsMeta.Range("A13").Value = $wsMeta.Range("A12").Value()
$wsMeta.Range("B13").Value = ""

# Row 12 now becomes the old row 11 (Description), but now with its
# description text filled in (it was blank before).
$wsMeta.Range("A12").Value = $wsMeta.Range("A11").Value()
$wsMeta.Range("B12").Value = "A set of codes further specifying the kind of Role; specific classification codes for further qualifying RoleClass codes."

# Row 11 becomes the brand new "Jurisdiction" row (value left blank).
$wsMeta.Range("A11").Value = "Jurisdiction"
$wsMeta.Range("B11").Value = ""

Write-Output "done"
